$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.194162845611572
$ws.Range("B1").Value = 2.342741012573242
$ws.Range("C1").Value = 3.569624185562134
$ws.Range("D1").Value = 3.229882955551147
$ws.Range("E1").Value = 1.138687133789062
